$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "56.534.06"
$ws.Range("E2").Value2 = "  +10.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.255.35"
$ws.Range("E3").Value2 = "  +6.48%  "
$ws.Range("E4").Value2 = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "399.28"
$ws.Range("E5").Value2 = "  +2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "111.35"
$ws.Range("E6").Value2 = "  +10.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.559"
$ws.Range("E7").Value2 = "  +4.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.999"
$ws.Range("E8").Value2 = "  -0.08%  "
$ws.Range("E9").Value2 = "  +7.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "39.54"
$ws.Range("E10").Value2 = "  +7.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0945"
$ws.Range("E11").Value2 = "  +11.72%  "
$ws.Range("E12").Value2 = "  +2.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "3.771.37"
$ws.Range("E13").Value2 = "  +6.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "19.25"
$ws.Range("E14").Value2 = "  +5.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "8.11"
$ws.Range("E15").Value2 = "  +6.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "3.255.16"
$ws.Range("E16").Value2 = "  +6.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "1.04"
$ws.Range("E17").Value2 = "  +5.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "10.96"
$ws.Range("E18").Value2 = "  +3.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "56.462.01"
$ws.Range("E19").Value2 = "  +10.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "3.32"
$ws.Range("E20").Value2 = "  +4.90%  "
$ws.Range("E21").Value2 = "  +8.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "13.05"
$ws.Range("E22").Value2 = "  +6.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "298.42"
$ws.Range("E23").Value2 = "  +13.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "75.00"
$ws.Range("E24").Value2 = "  +7.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "3.23"
$ws.Range("E25").Value2 = "  +2.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "8.10"
$ws.Range("E26").Value2 = "  +3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "28.22"
$ws.Range("E27").Value2 = "  +5.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "4.37"
$ws.Range("E28").Value2 = "  +5.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "7.34"
$ws.Range("E29").Value2 = "  +3.27%  "
$ws.Range("E30").Value2 = "  +4.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.00"
$ws.Range("E31").Value2 = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.112"
$ws.Range("E32").Value2 = "  +6.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "11.08"
$ws.Range("E33").Value2 = "  +6.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "38.33"
$ws.Range("E34").Value2 = "  +7.70%  "
$ws.Range("E35").Value2 = "  +0.06%  "
$ws.Range("E36").Value2 = "  +3.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "51.66"
$ws.Range("E37").Value2 = "  +3.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "3.14"
$ws.Range("E38").Value2 = "  +26.82%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value2 = "LidoDAOToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "3.51"
$ws.Range("E39").Value2 = "  +5.21%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value2 = "FirstDigitalUSD"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.00"
$ws.Range("E40").Value2 = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "17.59"
$ws.Range("E41").Value2 = "  +6.45%  "
$ws.Range("E42").Value2 = "  +6.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "133.52"
$ws.Range("E43").Value2 = "  +3.19%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value2 = "Stellar"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.120"
$ws.Range("E44").Value2 = "  +4.79%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value2 = "NEARProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "3.98"
$ws.Range("E45").Value2 = "  +5.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.286"
$ws.Range("E46").Value2 = "  -1.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "22.20"
$ws.Range("E47").Value2 = "  +2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "2.12"
$ws.Range("E48").Value2 = "  +49.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "2.152.54"
$ws.Range("E49").Value2 = "  +4.19%  "
$ws.Range("E50").Value2 = "  +1.36%  "
$ws.Range("E51").Value2 = "  -2.97%  "
